$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''27.435.77'
$ws.Range("E2").Value = '  +0.11%  '

# Row 3
$ws.Range("D3").Value = '''1.734.28'
$ws.Range("E3").Value = '  -0.83%  '

# Row 4
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '''322.92'
$ws.Range("E5").Value = '  +0.38%  '

# Row 6
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  -0.03%  '

# Row 7
$ws.Range("D7").Value = '''0.4530'
$ws.Range("E7").Value = '  +7.18%  '

# Row 8
$ws.Range("D8").Value = '''0.3525'
$ws.Range("E8").Value = '  -1.96%  '

# Row 9
$ws.Range("D9").Value = '''0.07389'
$ws.Range("E9").Value = '  -1.72%  '

# Row 10
$ws.Range("D10").Value = '''41.31'
$ws.Range("E10").Value = '  -2.53%  '

# Row 11
$ws.Range("E11").Value = '  -1.93%  '

# Row 12
$ws.Range("E12").Value = '  +0.00%  '

# Row 13
$ws.Range("D13").Value = '''20.31'
$ws.Range("E13").Value = '  -1.55%  '

# Row 14
$ws.Range("D14").Value = '''5.887'
$ws.Range("E14").Value = '  -2.38%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '''7.023'
$ws.Range("E15").Value = '  -2.46%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '''1.738.14'
$ws.Range("E16").Value = '  -0.91%  '

# Row 17
$ws.Range("D17").Value = '''91.20'
$ws.Range("E17").Value = '  -0.12%  '

# Row 18
$ws.Range("E18").Value = '  -1.23%  '

# Row 19
$ws.Range("D19").Value = '''0.06338'
$ws.Range("E19").Value = '  -0.21%  '

# Row 20
$ws.Range("D20").Value = '''1.002'
$ws.Range("E20").Value = '  -0.19%  '

# Row 21
$ws.Range("D21").Value = '''16.51'
$ws.Range("E21").Value = '  -2.98%  '

# Row 22
$ws.Range("D22").Value = '''5.703'
$ws.Range("E22").Value = '  -2.96%  '

# Row 23
$ws.Range("D23").Value = '''27.479.03'
$ws.Range("E23").Value = '  +0.18%  '

# Row 24
$ws.Range("D24").Value = '''11.06'
$ws.Range("E24").Value = '  -0.96%  '

# Row 25
$ws.Range("D25").Value = '''2.087'
$ws.Range("E25").Value = '  +0.00%  '

# Row 26
$ws.Range("D26").Value = '''161.70'
$ws.Range("E26").Value = '  +0.44%  '

# Row 27
$ws.Range("D27").Value = '''19.93'
$ws.Range("E27").Value = '  -1.44%  '

# Row 28
$ws.Range("D28").Value = '''1.940.69'
$ws.Range("E28").Value = '  -0.67%  '

# Row 29
$ws.Range("D29").Value = '''124.18'
$ws.Range("E29").Value = '  +0.78%  '

# Row 30
$ws.Range("D30").Value = '''2.025'
$ws.Range("E30").Value = '  -4.92%  '

# Row 31
$ws.Range("D31").Value = '''1.041'
$ws.Range("E31").Value = '  -6.01%  '

# Row 32
$ws.Range("D32").Value = '''0.09063'
$ws.Range("E32").Value = '  +2.64%  '

# Row 33
$ws.Range("D33").Value = '''3.652'
$ws.Range("E33").Value = '  +0.23%  '

# Row 34
$ws.Range("D34").Value = '''5.358'
$ws.Range("E34").Value = '  -3.24%  '

# Row 35
$ws.Range("D35").Value = '''0.02261'
$ws.Range("E35").Value = '  -0.59%  '

# Row 36
$ws.Range("D36").Value = '''11.56'
$ws.Range("E36").Value = '  -5.39%  '

# Row 37
$ws.Range("D37").Value = '''0.05932'
$ws.Range("E37").Value = '  -1.03%  '

# Row 38
$ws.Range("D38").Value = '''0.2050'
$ws.Range("E38").Value = '  -2.21%  '

# Row 39
$ws.Range("D39").Value = '''0.6205'
$ws.Range("E39").Value = '  -1.74%  '

# Row 40
$ws.Range("D40").Value = '''4.859'
$ws.Range("E40").Value = '  -1.42%  '

# Row 41
$ws.Range("D41").Value = '''1.186'
$ws.Range("E41").Value = '  +0.74%  '

# Row 42
$ws.Range("E42").Value = '  -1.12%  '

# Row 43
$ws.Range("D43").Value = '''7.662'
$ws.Range("E43").Value = '  -2.42%  '

# Row 44
$ws.Range("D44").Value = '''13.07'
$ws.Range("E44").Value = '  -2.42%  '

# Row 45
$ws.Range("D45").Value = '''3.700'
$ws.Range("E45").Value = '  +0.32%  '

# Row 46
$ws.Range("D46").Value = '''0.5763'
$ws.Range("E46").Value = '  -1.46%  '

# Row 47
$ws.Range("D47").Value = '''121.93'
$ws.Range("E47").Value = '  -0.32%  '

# Row 48
$ws.Range("D48").Value = '''1.924'
$ws.Range("E48").Value = '  -2.20%  '

# Row 49
$ws.Range("E49").Value = '  +0.70%  '

# Row 50
$ws.Range("D50").Value = '''1.106'
$ws.Range("E50").Value = '  -3.86%  '

# Row 51
$ws.Range("D51").Value = '''70.95'
$ws.Range("E51").Value = '  -2.84%  '
